# Update the "Score" column (C) values on Sheet1 for the Nursing VR Study
# workbook to reflect the corrected/updated midterm exam data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$scores = @{
    2  = 26
    3  = 31
    4  = 35
    5  = 37
    6  = 30
    7  = 29
    8  = 25
    9  = 26
    10 = 28
    11 = 23
    12 = 35
    13 = 27
    14 = 41
    15 = 38
    16 = 29
    17 = 32
    18 = 45
    19 = 33
    20 = 15
    22 = 37
    23 = 42
    24 = 37
    26 = 47
    27 = 31
    28 = 29
    29 = 26
    30 = 40
    31 = 27
    32 = 49
    33 = 39
    35 = 46
    36 = 33
    37 = 44
    38 = 42
    39 = 30
    40 = 40
    41 = 45
}

foreach ($row in $scores.Keys) {
    $ws.Cells.Item($row, 3).Value = $scores[$row]
}
